$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New "Praktiken" rows (3-11): label in column A, "X" marks in B/C/D
# ---------------------------------------------------------------------------
$ws.Range("A3").Value()  = "Informations-Sicherheit"
$ws.Range("A4").Value()  = "Wissensmanagement"
$ws.Range("A5").Value()  = "Messen und Dokumentieren"
$ws.Range("A6").Value()  = "Risikomanagement"
$ws.Range("A7").Value()  = "Zulieferermanagement"
$ws.Range("A8").Value()  = "Incidentmanagement"
$ws.Range("A9").Value()  = "IT-Kapitalmanagement"
$ws.Range("A10").Value() = "Problemmanagement"
$ws.Range("A11").Value() = "Releasemanagement"

$ws.Range("B3:D5").Value()  = "X"
$ws.Range("C6:D7").Value()  = "X"
$ws.Range("B8:D8").Value()  = "X"
$ws.Range("C9:D9").Value()  = "X"
$ws.Range("B10").Value()    = "X"
$ws.Range("D10").Value()    = "X"
$ws.Range("D11").Value()    = "X"

# Center the new "X" marker cells, like the rest of the matrix
$ws.Range("B3:D5").HorizontalAlignment()  = -4108
$ws.Range("B3:D5").VerticalAlignment()    = -4108
$ws.Range("C6:D7").HorizontalAlignment()  = -4108
$ws.Range("C6:D7").VerticalAlignment()    = -4108
$ws.Range("B8:D8").HorizontalAlignment()  = -4108
$ws.Range("B8:D8").VerticalAlignment()    = -4108
$ws.Range("C9:D9").HorizontalAlignment()  = -4108
$ws.Range("C9:D9").VerticalAlignment()    = -4108
$ws.Range("B10").HorizontalAlignment()    = -4108
$ws.Range("B10").VerticalAlignment()      = -4108
$ws.Range("D10").HorizontalAlignment()    = -4108
$ws.Range("D10").VerticalAlignment()      = -4108
$ws.Range("D11").HorizontalAlignment()    = -4108
$ws.Range("D11").VerticalAlignment()      = -4108

# ---------------------------------------------------------------------------
# Remove the stray header fill bleeding past column D (was filled purple/pink,
# should be unfilled like the rest of the sheet)
# ---------------------------------------------------------------------------
$ws.Range("E1:M2").Interior.Pattern() = -4142

# ---------------------------------------------------------------------------
# Column widths: widen the label column, set a sensible width for the new
# X/O columns
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth() = 24.1667
$ws.Range("B1:D1").EntireColumn.ColumnWidth() = 10.6667

# ---------------------------------------------------------------------------
# Selection cosmetics (matches the author's final cursor position)
# ---------------------------------------------------------------------------
$ws.Range("E11").Select()
